$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------
# 1) Window view tweak: windowHeight 15720 -> 15840 is not exposed
#    via the object model in this runtime (cosmetic window chrome
#    size); skipped as not controllable.
# ---------------------------------------------------------------

# ---------------------------------------------------------------
# 2) Row 2 loses its explicit 30pt height (goes back to auto height)
# ---------------------------------------------------------------
$ws.Rows(2).AutoFit()

# ---------------------------------------------------------------
# 3) Row 29: height 60 -> 45 (content/date/hours stay the same)
# ---------------------------------------------------------------
$ws.Rows(29).RowHeight = 45

# ---------------------------------------------------------------
# 4) Insert a brand-new row after row 29 (becomes the new row 30)
#    for the new task "Regler probleme crash application".
#    Inserting copies formatting from the row above (row 29).
# ---------------------------------------------------------------
$ws.Rows(30).Insert()
$ws.Range("B30").Value = "3h"
$ws.Range("C30").Value = "Régler problème crash application"

# The old row 30 (now row 31) keeps ht=30 already (inherited); nothing to do.
# The old row 31 (now row 32) needs its height changed from 30 -> 45.
$ws.Rows(32).RowHeight = 45

# ---------------------------------------------------------------
# 5) Insert two brand-new rows after row 32 (new rows 33 and 34)
#    for the two additional tasks.
# ---------------------------------------------------------------
$ws.Rows(33).Insert()
$ws.Range("B33").Value = "2h"
$ws.Range("C33").Value = "Implémentation menus déroulants de sélection et de tri des indices"
$ws.Rows(33).RowHeight = 30

$ws.Rows(34).Insert()
$ws.Range("B34").Value = "1h"
$ws.Range("C34").Value = "Continuation du guide d'installation"

# ---------------------------------------------------------------
# 6) Merge cells: the date cell for the new row 30 merges with A29,
#    and the date cell area for rows 31-34 merges together.
# ---------------------------------------------------------------
$ws.Range("A31:A32").UnMerge()
$ws.Range("A29:A30").Merge()
$ws.Range("A31:A34").Merge()

# ---------------------------------------------------------------
# 7) Selection / scroll position shown in the file when last saved.
# ---------------------------------------------------------------
$ws.Range("D29").Select()

Write-Output "edit applied"
